$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H20").Value = 883
$ws.Range("I20").Value = 510.66666
$ws.Range("J20").Value = 2000
$ws.Range("K20").Value = 510.66666
$ws.Range("L20").Value = 2000
$ws.Range("M20").Value = -280.66666
$ws.Range("N20").Value = -2460
$ws.Range("H35").Value = 883
$ws.Range("I35").Value = 510.66666
$ws.Range("J35").Value = 2000
$ws.Range("K35").Value = 510.66666
$ws.Range("L35").Value = 2000
$ws.Range("M35").Value = -131.66666
$ws.Range("N35").Value = -2758
$ws.Range("H137").Value = 1501.7646
$ws.Range("J137").Value = 1493.3334
$ws.Range("L137").Value = 4480.0002
$ws.Range("N137").Value = -9580.0002
$ws.Range("H138").Value = 3847.5652
$ws.Range("I138").Value = 3149.3333
$ws.Range("J138").Value = 4296.4287
$ws.Range("K138").Value = 9447.999899999999
$ws.Range("L138").Value = 12889.2861
$ws.Range("M138").Value = -4307.999899999999
$ws.Range("N138").Value = -23169.2861

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3568.8
$ws.Range("I32").Value = 2454.261
$ws.Range("J32").Value = 7230.857
$ws.Range("K32").Value = 2454.261
$ws.Range("L32").Value = 7230.857
$ws.Range("M32").Value = -2167.261
$ws.Range("N32").Value = -7804.857
$ws.Range("H74").Value = 1000
$ws.Range("I74").Value = 1000
$ws.Range("K74").Value = 1000
$ws.Range("M74").Value = -126
$ws.Range("H77").Value = 1000
$ws.Range("I77").Value = 1000
$ws.Range("K77").Value = 5000
$ws.Range("M77").Value = -632

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 7998.6
$ws.Range("I20").Value = 6249.5
$ws.Range("J20").Value = 14995
$ws.Range("K20").Value = 6249.5
$ws.Range("L20").Value = 14995
$ws.Range("M20").Value = -6002.5
$ws.Range("N20").Value = -15489
$ws.Range("H80").Value = 482.7143
$ws.Range("I80").Value = 526
$ws.Range("K80").Value = 526
$ws.Range("M80").Value = 472
$ws.Range("H83").Value = 482.7143
$ws.Range("I83").Value = 526
$ws.Range("K83").Value = 2630
$ws.Range("M83").Value = 2362
$ws.Range("H105").Value = 3530.4443
$ws.Range("I105").Value = 3471.75
$ws.Range("J105").Value = 4000
$ws.Range("K105").Value = 3471.75
$ws.Range("L105").Value = 4000
$ws.Range("M105").Value = -1724.75
$ws.Range("N105").Value = -7494

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3166.6667
$ws.Range("I31").Value = 1500
$ws.Range("J31").Value = 4000
$ws.Range("K31").Value = 1500
$ws.Range("L31").Value = 4000
$ws.Range("M31").Value = -1205
$ws.Range("N31").Value = -4590
$ws.Range("H34").Value = 3166.6667
$ws.Range("I34").Value = 1500
$ws.Range("J34").Value = 4000
$ws.Range("K34").Value = 1500
$ws.Range("L34").Value = 4000
$ws.Range("M34").Value = -1298
$ws.Range("N34").Value = -4404
$ws.Range("H132").Value = 2003.6072
$ws.Range("I132").Value = 2011.1852
$ws.Range("J132").Value = 1799
$ws.Range("K132").Value = 6033.5556
$ws.Range("L132").Value = 5397
$ws.Range("M132").Value = -3503.5556
$ws.Range("N132").Value = -10457
$ws.Range("H134").Value = 1658
$ws.Range("I134").Value = 1435.2
$ws.Range("K134").Value = 4305.6
$ws.Range("M134").Value = -1770.6

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 284.5
$ws.Range("J12").Value = 284.5
$ws.Range("L12").Value = 853.5
$ws.Range("N12").Value = -1199.5
$ws.Range("H107").Value = 0
$ws.Range("I107").Value = 0
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 0
$ws.Range("L107").ClearContents()
$ws.Range("M107").ClearContents()
$ws.Range("N107").Value = 0
$ws.Range("H124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("L124").ClearContents()
$ws.Range("N124").Value = 0
$ws.Range("H125").Value = 25000
$ws.Range("J125").Value = 25000
$ws.Range("L125").Value = 75000
$ws.Range("N125").Value = -84840
$ws.Range("H126").Value = 2500
$ws.Range("I126").Value = 2500
$ws.Range("K126").Value = 7500
$ws.Range("M126").Value = -2560
$ws.Range("H129").Value = 2357
$ws.Range("J129").Value = 2697.3333
$ws.Range("L129").Value = 8091.999899999999
$ws.Range("N129").Value = -18091.9999
$ws.Range("H130").Value = 14500
$ws.Range("J130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("N130").ClearContents()
$ws.Range("H131").Value = 873
$ws.Range("I131").Value = 873
$ws.Range("J131").Value = 0
$ws.Range("K131").Value = 2619
$ws.Range("L131").Value = 0
$ws.Range("M131").ClearContents()
$ws.Range("N131").Value = 2421
$ws.Range("H139").Value = 2364.3333
$ws.Range("I139").Value = 2152.4443
$ws.Range("K139").Value = 6457.3329
$ws.Range("M139").Value = -1317.3329

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 10817.875
$ws.Range("I70").Value = 10963.286
$ws.Range("J70").Value = 9800
$ws.Range("K70").Value = 10963.286
$ws.Range("L70").Value = 9800
$ws.Range("M70").Value = -10693.286
$ws.Range("N70").Value = -10340
$ws.Range("H73").Value = 10817.875
$ws.Range("I73").Value = 10963.286
$ws.Range("J73").Value = 9800
$ws.Range("K73").Value = 10963.286
$ws.Range("L73").Value = 9800
$ws.Range("M73").Value = -10027.286
$ws.Range("N73").Value = -11672
$ws.Range("H132").Value = 2209.7273
$ws.Range("I132").Value = 2339.2
$ws.Range("K132").Value = 7017.599999999999
$ws.Range("M132").Value = -4487.599999999999

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 890.125
$ws.Range("J55").Value = 913.8570999999999
$ws.Range("L55").Value = 913.8570999999999
$ws.Range("N55").Value = -1259.8571
$ws.Range("H132").Value = 475
$ws.Range("I132").Value = 475
$ws.Range("K132").Value = 1425
$ws.Range("M132").Value = 1105

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H119").Value = 52739.4
$ws.Range("J119").Value = 52739.4
$ws.Range("L119").Value = 52739.4
$ws.Range("N119").Value = -62415.4
$ws.Range("H132").Value = 2403.7693
$ws.Range("I132").Value = 2154.9
$ws.Range("K132").Value = 6464.700000000001
$ws.Range("M132").Value = -3934.700000000001
